$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 2288
$ws.Range("J3").Value = 8078
$ws.Range("K3").Value = 2200
$ws.Range("B4").Value = 1698
$ws.Range("K4").Value = 469
$ws.Range("J5").Value = 626
$ws.Range("K6").Value = 2768
$ws.Range("B7").Value = 23331
$ws.Range("K7").Value = 7871

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 58
$ws.Range("J6").Value = 228
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 230
$ws.Range("J8").Value = 1853
$ws.Range("K8").Value = 528
$ws.Range("K9").Value = 31
$ws.Range("K10").Value = 44
$ws.Range("K11").Value = 168
$ws.Range("K14").Value = 44
$ws.Range("J15").Value = 358
$ws.Range("K15").Value = 78
$ws.Range("K16").Value = 18
$ws.Range("K19").Value = 226
$ws.Range("K20").Value = 171
$ws.Range("K23").Value = 71
$ws.Range("K29").Value = 398
$ws.Range("K31").Value = 87
$ws.Range("K33").Value = 308
$ws.Range("K37").Value = 253
$ws.Range("K42").Value = 271
$ws.Range("K51").Value = 86
$ws.Range("K52").Value = 212
$ws.Range("K53").Value = 114
$ws.Range("K54").Value = 148
$ws.Range("K55").Value = 87
$ws.Range("K57").Value = 23
$ws.Range("K60").Value = 54
$ws.Range("B63").Value = 403
$ws.Range("J63").Value = 98
$ws.Range("K63").Value = 33
$ws.Range("K67").Value = 303
$ws.Range("K68").Value = 20
$ws.Range("K74").Value = 6
$ws.Range("K81").Value = 8
$ws.Range("K83").Value = 170
$ws.Range("K84").Value = 55
$ws.Range("K85").Value = 387
$ws.Range("K86").Value = 55
$ws.Range("K88").Value = 98
$ws.Range("K89").Value = 106
$ws.Range("K91").Value = 73
$ws.Range("K96").Value = 107
$ws.Range("K99").Value = 142
$ws.Range("B101").Value = 23331
$ws.Range("K101").Value = 7871

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 72
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 52
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 140
$ws.Range("K3").Value = 130
$ws.Range("K7").Value = 387

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 156
$ws.Range("K3").Value = 153
$ws.Range("J4").Value = 97
$ws.Range("K4").Value = 28
$ws.Range("K6").Value = 179
$ws.Range("J7").Value = 1853
$ws.Range("K7").Value = 528

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 88
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 308

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 35
$ws.Range("K4").Value = 7

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 64
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 303

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 27
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 106
$ws.Range("K3").Value = 131
$ws.Range("K7").Value = 398

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 64
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 226

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K3").Value = 20
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 228
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 271

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 31
$ws.Range("K6").Value = 38

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 34
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 24
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 22
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 73
$ws.Range("K6").Value = 27
$ws.Range("J7").Value = 358
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 8
